$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp label (row 1)
$ws.Range("A1").Value = "Datos actualizados a 7 de Agosto de 2020 a las 06:17"

# Update country rows: name (column A) plus Casos totales/Nuevos casos/Casos
# activos/Recuperados/Casos criticos/Muertes hoy/Muertes (columns B:H)
$updates = @(
    @{ Row = 29; Name = "Kazajistan"; Vals = @(96922, 980, 70680, 25184, 0, 0, 1058) },
    @{ Row = 32; Name = "Bolivia"; Vals = @(86423, 1282, 27373, 55585, 0, 80, 3465) },
    @{ Row = 39; Name = "Belgica"; Vals = @(72016, 858, 17700, 44455, 0, 2, 9861) },
    @{ Row = 40; Name = "Panama"; Vals = @(71418, 0, 45658, 24186, 0, 0, 1574) },
    @{ Row = 50; Name = "Honduras"; Vals = @(45755, 657, 6225, 38084, 0, 23, 1446) },
    @{ Row = 51; Name = "Nigeria"; Vals = @(45244, 0, 32430, 11884, 0, 0, 930) },
    @{ Row = 67; Name = "Venezuela"; Vals = @(23280, 0, 12470, 10608, 0, 0, 202) },
    @{ Row = 72; Name = "Australia"; Vals = @(20270, 408, 11147, 8857, 0, 11, 266) },
    @{ Row = 93; Name = "Haiti"; Vals = @(7582, 38, 4832, 2579, 0, 0, 171) },
    @{ Row = 159; Name = "Vietnam"; Vals = @(750, 3, 392, 348, 0, 0, 10) },
    @{ Row = 167; Name = "Burundi"; Vals = @(400, 0, 304, 95, 0, 0, 1) },
    @{ Row = 168; Name = "Comoras"; Vals = @(396, 0, 340, 49, 0, 0, 7) },
    @{ Row = 182; Name = "San Martin (Parte Holandesa)"; Vals = @(176, 16, 86, 74, 0, 0, 16) },
    @{ Row = 183; Name = "Papua Nueva Guinea"; Vals = @(163, 0, 53, 107, 0, 0, 3) },
    @{ Row = 185; Name = "Islas Turcas y Caicos"; Vals = @(141, 12, 39, 100, 0, 0, 2) },
    @{ Row = 186; Name = "Brunei"; Vals = @(141, 0, 138, 0, 0, 0, 3) },
    @{ Row = 187; Name = "Barbados"; Vals = @(133, 0, 100, 26, 0, 0, 7) },
    @{ Row = 190; Name = "Belice"; Vals = @(114, 28, 32, 80, 0, 0, 2) },
    @{ Row = 191; Name = "Butan"; Vals = @(108, 3, 96, 12, 0, 0, 0) },
    @{ Row = 192; Name = "Antigua y Barbuda"; Vals = @(92, 0, 76, 13, 0, 0, 3) },
    @{ Row = 193; Name = "Liechtenstein"; Vals = @(89, 0, 85, 3, 0, 0, 1) },
    @{ Row = 202; Name = "Santa Lucia"; Vals = @(25, 0, 24, 1, 0, 0, 0) },
    @{ Row = 203; Name = "Timor Oriental"; Vals = @(25, 0, 24, 1, 0, 0, 0) }
)

foreach ($u in $updates) {
    $r = $u.Row
    $ws.Range("A$r").Value = $u.Name
    $vals = $u.Vals
    for ($i = 0; $i -lt $vals.Length; $i++) {
        $ws.Cells.Item($r, $i + 2).Value = $vals[$i]
    }
}
